$d = $word.ActiveDocument

$replacements = @(
    @("266×7=1862", "642×3=1926"),
    @("423×5=2115", "643×4=2572"),
    @("927×4=3708", "683×8=5464"),
    @("729×8=5832", "392×5=1960"),
    @("788×9=7092", "917×3=2751"),
    @("561×6=3366", "773×8=6184"),
    @("820×5=4100", "298×8=2384"),
    @("105×8=840",  "137×6=822"),
    @("105×2=210",  "354×2=708"),
    @("386×5=1930", "565×2=1130"),
    @("933×4=3732", "841×4=3364"),
    @("297×3=891",  "347×6=2082"),
    @("870×2=1740", "899×8=7192"),
    @("947×2=1894", "298×4=1192"),
    @("341×7=2387", "369×5=1845"),
    @("444×9=3996", "393×8=3144"),
    @("629×3=1887", "792×9=7128"),
    @("755×8=6040", "434×4=1736"),
    @("688×3=2064", "481×8=3848"),
    @("583×2=1166", "763×6=4578"),
    @("863×2=1726", "508×4=2032"),
    @("713×5=3565", "436×8=3488"),
    @("192×7=1344", "142×5=710"),
    @("275×3=825",  "590×7=4130"),
    @("531×9=4779", "290×6=1740")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
